$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.ClearFormats()
}

Set-TextValue "D2" "56.544.94"
Set-TextValue "E2" "  +4.08%  "
Set-TextValue "D3" "3.001.58"
Set-TextValue "E3" "  +4.51%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "507.60"
Set-TextValue "E5" "  +8.09%  "
Set-TextValue "D6" "136.76"
Set-TextValue "E6" "  +8.54%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "0.433"
Set-TextValue "E8" "  +7.40%  "
Set-TextValue "D9" "7.59"
Set-TextValue "E9" "  +14.61%  "
Set-TextValue "E10" "  +12.62%  "
Set-TextValue "D11" "0.353"
Set-TextValue "E11" "  +7.57%  "
Set-TextValue "E12" "  +5.28%  "
Set-TextValue "D13" "3.516.29"
Set-TextValue "E13" "  +4.64%  "
Set-TextValue "D14" "25.51"
Set-TextValue "E14" "  +10.31%  "
Set-TextValue "E15" "  +15.79%  "
Set-TextValue "D16" "56.564.33"
Set-TextValue "E16" "  +4.18%  "
Set-TextValue "D17" "3.000.33"
Set-TextValue "E17" "  +4.31%  "
Set-TextValue "D18" "5.86"
Set-TextValue "E18" "  +10.01%  "
Set-TextValue "D19" "12.53"
Set-TextValue "E19" "  +10.00%  "
Set-TextValue "D20" "7.82"
Set-TextValue "E20" "  +10.98%  "
Set-TextValue "D21" "326.77"
Set-TextValue "E21" "  +10.00%  "
Set-TextValue "E22" "  -0.05%  "
Set-TextValue "D23" "0.478"
Set-TextValue "E23" "  +8.71%  "
Set-TextValue "D24" "62.48"
Set-TextValue "E24" "  +6.86%  "
Set-TextValue "D25" "0.167"
Set-TextValue "E25" "  +10.00%  "
Set-TextValue "E26" "  +0.06%  "
Set-TextValue "D27" "0.0₃0914"
Set-TextValue "E27" "  +13.68%  "
Set-TextValue "D28" "6.54"
Set-TextValue "E28" "  +7.25%  "
Set-TextValue "D29" "6.96"
Set-TextValue "D30" "1.24"
Set-TextValue "E30" "  +11.99%  "
Set-TextValue "E31" "  +9.70%  "
Set-TextValue "D32" "20.66"
Set-TextValue "E32" "  +11.30%  "
Set-TextValue "D33" "155.34"
Set-TextValue "E33" "  +13.72%  "
Set-TextValue "D34" "4.50"
Set-TextValue "E34" "  +8.16%  "
Set-TextValue "D35" "5.63"
Set-TextValue "E35" "  +4.78%  "
Set-TextValue "D36" "1.27"
Set-TextValue "E36" "  +4.36%  "
Set-TextValue "D37" "0.0672"
Set-TextValue "E37" "  +9.65%  "
Set-TextValue "D38" "23.96"
Set-TextValue "E38" "  +4.39%  "
Set-TextValue "D39" "3.038.45"
Set-TextValue "E39" "  +4.96%  "
Set-TextValue "D40" "36.63"
Set-TextValue "E40" "  +5.01%  "
Set-TextValue "E41" "  -0.07%  "
Set-TextValue "E42" "  +7.98%  "
Set-TextValue "D43" "2.265.16"
Set-TextValue "E43" "  +11.11%  "
Set-TextValue "D44" "0.996"
Set-TextValue "E44" "  +5.73%  "
Set-TextValue "E45" "  +7.38%  "
Set-TextValue "D46" "3.61"
Set-TextValue "E46" "  +7.41%  "
Set-TextValue "D47" "1.97"
Set-TextValue "E47" "  +23.10%  "
Set-TextValue "E48" "  +10.91%  "
Set-TextValue "D49" "5.78"
Set-TextValue "E49" "  +8.24%  "
Set-TextValue "D50" "19.11"
Set-TextValue "E50" "  +7.54%  "
Set-TextValue "E51" "  +11.79%  "
